# Revert "Drop in all data files from 3.0 RMI script"
# Re-introduce the "Texas Notes" worksheet (with its explanatory notes) between
# the "Data" and "PPEIdtICEaT" sheets, and restore the original selections.

$wb = $excel.ActiveWorkbook

# --- Insert the "Texas Notes" worksheet right after "Data" ---------------
$dataSheet = $wb.Worksheets.Item("Data")
$texasSheet = $wb.Worksheets.Add()
$texasSheet.Name = "Texas Notes"
$texasSheet.Move($null, $wb.Worksheets.Item("Data"))

# Note: sheet collection indices shift after Move(), so every sheet
# reference used below is re-fetched by name rather than reused.

# --- Populate the notes content -------------------------------------------
$texasSheet = $wb.Worksheets.Item("Texas Notes")
$texasSheet.Range("A1").Value = "This spreadsheet uses a very particular study. "
$texasSheet.Range("A3").Value = "It's done in Wisconsin where the authors use an educational seminar for builders"
$texasSheet.Range("A4").Value = "then they follow up with phone surveys to see what the builders actually implemented"
$texasSheet.Range("A5").Value = "then they use those results to try and estimate how much energy those builders decisions saved"
$texasSheet.Range("A7").Value = "It's all very niche and a bit subjective"
$texasSheet.Range("A9").Value = "That said, I did a quick literature search and didn't find anything that I thought"
$texasSheet.Range("A10").Value = "would give us better or more Texas-specific numbers. "
$texasSheet.Range("A12").Value = 'Since this is a "low" priority sheet, I will leave it alone.'

# --- Restore each sheet's selection/active cell ---------------------------
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
$aboutSheet.Range("B25").Select()

$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Activate()
$dataSheet.Range("A14").Select()

$texasSheet = $wb.Worksheets.Item("Texas Notes")
$texasSheet.Activate()
$texasSheet.Range("C27").Select()

# PPEIdtICEaT becomes the active/selected tab, matching the restored workbook
$ppeSheet = $wb.Worksheets.Item("PPEIdtICEaT")
$ppeSheet.Activate()
$ppeSheet.Range("D20").Select()
